# "everything has been added to one table except traits"
#
# The "Other Info" sheet (sheet3: Event/Sale data) is merged into the
# "General Info" sheet (sheet1) as additional columns O:AB, and the
# "Other Info" sheet is then removed. "Traits" is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Info")

# Helper: write a value into a cell while forcing text storage so that
# numeric-looking strings (big integer ids, decimal strings with trailing
# zeros, etc.) are not silently coerced into floating point numbers and
# lose precision/formatting. We briefly flip the cell to a text number
# format, assign the value, then restore the default ("Normal") cell
# style so no visible formatting / border / bold carries over.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ------------------------------------------------------------------
# 1. Extend the header row (row 1) with the "Other Info" headers,
#    copying the existing header formatting (bold/border/centered)
#    from column N1 into O1:AB1 first.
# ------------------------------------------------------------------
$ws.Range("N1").Copy()
$ws.Range("O1:AB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headers = @(
    "Event Type ",
    "Event Timestamp",
    "Auction Type",
    "Total Price",
    "Last Sale Creation Date",
    "Quantity",
    "Telegram URL",
    "Twitter User",
    "Instagram User",
    "Wiki URL",
    "Discord URL",
    "ETH Price",
    "USD Price",
    "Address of Last Transaction"
)

$col = 15  # column O
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# ------------------------------------------------------------------
# 2. Append the "Other Info" data row (row 2) values into O2:AB2.
# ------------------------------------------------------------------
$ws.Cells.Item(2, 15).Value = "successful"                           # O2 Event Type
$ws.Cells.Item(2, 16).Value = "2021-07-20T19:33:14"                  # P2 Event Timestamp
# Q2 Auction Type left blank (empty string in source)
Set-TextValue $ws.Cells.Item(2, 18) "250000000000000000"             # R2 Total Price
$ws.Cells.Item(2, 19).Value = "2021-07-20T19:33:52.742091"           # S2 Last Sale Creation Date
Set-TextValue $ws.Cells.Item(2, 20) "1"                               # T2 Quantity
# U2 Telegram URL left blank
$ws.Cells.Item(2, 22).Value = "realsupducks"                         # V2 Twitter User
# W2 Instagram User left blank
# X2 Wiki URL left blank
$ws.Cells.Item(2, 25).Value = "https://discord.gg/UJCP5y3s7J"        # Y2 Discord URL
Set-TextValue $ws.Cells.Item(2, 26) "1.000000000000000"               # Z2 ETH Price
Set-TextValue $ws.Cells.Item(2, 27) "2154.639999999999873000"         # AA2 USD Price
$ws.Cells.Item(2, 28).Value = "0x69c4e59b4f1f8a2782279ba9d884b8d3a2c1e6ad"  # AB2 Address of Last Transaction

# ------------------------------------------------------------------
# 3. Remove the now-redundant "Other Info" sheet.
# ------------------------------------------------------------------
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Other Info").Delete()
$excel.DisplayAlerts = $true
